$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1575.2222
$ws.Range("J43").Value = 1961.75
$ws.Range("L43").Value = 1961.75
$ws.Range("N43").Value = -2099.75
$ws.Range("H58").Value = 380
$ws.Range("I58").Value = 380
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 1140
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -990
$ws.Range("N58").ClearContents()
$ws.Range("H69").Value = 7984.909
$ws.Range("H72").Value = 7984.909
$ws.Range("H76").Value = 2319970.2
$ws.Range("I76").Value = 3476612.5
$ws.Range("K76").Value = 3476612.5
$ws.Range("M76").Value = -3476297.5
$ws.Range("H79").Value = 2319970.2
$ws.Range("I79").Value = 3476612.5
$ws.Range("K79").Value = 3476612.5
$ws.Range("M79").Value = -3475520.5
$ws.Range("H100").Value = 1922.5
$ws.Range("I100").Value = 1734.6364
$ws.Range("K100").Value = 1734.6364
$ws.Range("M100").Value = -1193.6364
$ws.Range("H106").Value = 83335330
$ws.Range("I106").Value = 100001400
$ws.Range("K106").Value = 100001400
$ws.Range("M106").Value = -100000769
$ws.Range("H138").Value = 7839.8706
$ws.Range("J138").Value = 7999.1567
$ws.Range("L138").Value = 23997.4701
$ws.Range("N138").Value = -34277.4701
$ws.Range("H140").Value = 59987
$ws.Range("J140").Value = 59987
$ws.Range("L140").Value = 59987
$ws.Range("N140").Value = -70347

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1489703.1
$ws.Range("I2").Value = 1951393.5
$ws.Range("J2").Value = 2034.3334
$ws.Range("K2").Value = 1951393.5
$ws.Range("L2").Value = 2034.3334
$ws.Range("M2").Value = -1951280.5
$ws.Range("N2").Value = -2260.3334
$ws.Range("H32").Value = 16395.371
$ws.Range("I32").Value = 13145.841
$ws.Range("K32").Value = 13145.841
$ws.Range("M32").Value = -12858.841
$ws.Range("H45").Value = 5716825
$ws.Range("J45").Value = 7994.6
$ws.Range("L45").Value = 7994.6
$ws.Range("N45").Value = -8748.6
$ws.Range("H46").Value = 10449.667
$ws.Range("I46").Value = 8625
$ws.Range("J46").Value = 11362
$ws.Range("K46").Value = 8625
$ws.Range("L46").Value = 11362
$ws.Range("M46").Value = -8306
$ws.Range("N46").Value = -12000
$ws.Range("H74").Value = 85386.73
$ws.Range("I74").Value = 5379.778
$ws.Range("J74").Value = 445418
$ws.Range("K74").Value = 5379.778
$ws.Range("L74").Value = 445418
$ws.Range("M74").Value = -4505.778
$ws.Range("N74").Value = -447166
$ws.Range("H77").Value = 85386.73
$ws.Range("I77").Value = 5379.778
$ws.Range("J77").Value = 445418
$ws.Range("K77").Value = 26898.89
$ws.Range("L77").Value = 2227090
$ws.Range("M77").Value = -22530.89
$ws.Range("N77").Value = -2235826
$ws.Range("H116").Value = 1489703.1
$ws.Range("I116").Value = 1951393.5
$ws.Range("J116").Value = 2034.3334
$ws.Range("K116").Value = 1951393.5
$ws.Range("L116").Value = 2034.3334
$ws.Range("M116").Value = -1949099.5
$ws.Range("N116").Value = -6622.3334
$ws.Range("H132").Value = 35869.367
$ws.Range("I132").Value = 2746.2666
$ws.Range("K132").Value = 8238.799800000001
$ws.Range("M132").Value = -5708.799800000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1489703.1
$ws.Range("I3").Value = 1951393.5
$ws.Range("J3").Value = 2034.3334
$ws.Range("K3").Value = 1951393.5
$ws.Range("L3").Value = 2034.3334
$ws.Range("M3").Value = -1951279.5
$ws.Range("N3").Value = -2262.3334
$ws.Range("H10").Value = 2478
$ws.Range("J10").Value = 2478
$ws.Range("L10").Value = 2478
$ws.Range("N10").Value = -2758
$ws.Range("H99").Value = 11989680
$ws.Range("I99").Value = 15985402
$ws.Range("J99").Value = 2511.3333
$ws.Range("K99").Value = 15985402
$ws.Range("L99").Value = 2511.3333
$ws.Range("M99").Value = -15983904
$ws.Range("N99").Value = -5507.3333

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 24601.373
$ws.Range("I31").Value = 3396.1538
$ws.Range("J31").Value = 33790.3
$ws.Range("K31").Value = 3396.1538
$ws.Range("L31").Value = 33790.3
$ws.Range("M31").Value = -3101.1538
$ws.Range("N31").Value = -34380.3
$ws.Range("H34").Value = 24601.373
$ws.Range("I34").Value = 3396.1538
$ws.Range("J34").Value = 33790.3
$ws.Range("K34").Value = 3396.1538
$ws.Range("L34").Value = 33790.3
$ws.Range("M34").Value = -3194.1538
$ws.Range("N34").Value = -34194.3
$ws.Range("H132").Value = 101922.48
$ws.Range("I132").Value = 68875.734
$ws.Range("J132").Value = 225847.75
$ws.Range("K132").Value = 206627.202
$ws.Range("L132").Value = 677543.25
$ws.Range("M132").Value = -204097.202
$ws.Range("N132").Value = -682603.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 127003.125
$ws.Range("J55").Value = 253500
$ws.Range("L55").Value = 760500
$ws.Range("N55").Value = -760854
$ws.Range("H122").Value = 2134.8
$ws.Range("J122").Value = 2889.8
$ws.Range("L122").Value = 26008.2
$ws.Range("N122").Value = -30908.2
$ws.Range("H137").Value = 4444.6113
$ws.Range("J137").Value = 5498.0835
$ws.Range("L137").Value = 16494.2505
$ws.Range("N137").Value = -26694.2505

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H54").Value = 12861.857
$ws.Range("J54").Value = 12861.857
$ws.Range("L54").Value = 12861.857
$ws.Range("N54").Value = -13641.857
$ws.Range("H102").Value = 4998751
$ws.Range("I102").Value = 5848724
$ws.Range("J102").Value = 2980064.5
$ws.Range("K102").Value = 5848724
$ws.Range("L102").Value = 2980064.5
$ws.Range("M102").Value = -5847102
$ws.Range("N102").Value = -2983308.5
$ws.Range("H113").Value = 7625558
$ws.Range("I113").Value = 15244127
$ws.Range("J113").Value = 6989.909
$ws.Range("K113").Value = 15244127
$ws.Range("L113").Value = 6989.909
$ws.Range("M113").Value = -15241957
$ws.Range("N113").Value = -11329.909
$ws.Range("H126").Value = 4956515.5
$ws.Range("J126").Value = 7579594
$ws.Range("L126").Value = 22738782
$ws.Range("N126").Value = -22743722

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10742.833
$ws.Range("J7").Value = 11136.857
$ws.Range("L7").Value = 11136.857
$ws.Range("N7").Value = -11360.857
$ws.Range("H40").Value = 7975.769
$ws.Range("I40").Value = 5148.3076
$ws.Range("J40").Value = 10803.23
$ws.Range("K40").Value = 5148.3076
$ws.Range("L40").Value = 10803.23
$ws.Range("M40").Value = -5012.3076
$ws.Range("N40").Value = -11075.23
$ws.Range("H100").Value = 145386.28
$ws.Range("I100").Value = 3160
$ws.Range("J100").Value = 500952
$ws.Range("K100").Value = 3160
$ws.Range("L100").Value = 500952
$ws.Range("M100").Value = -2619
$ws.Range("N100").Value = -502034
$ws.Range("H126").Value = 10742.833
$ws.Range("J126").Value = 11136.857
$ws.Range("L126").Value = 33410.571
$ws.Range("N126").Value = -38350.571

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7939774
$ws.Range("I81").Value = 7939774
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 15879548
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -15878487
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 7939774
$ws.Range("I84").Value = 7939774
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 79397740
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -79392436
$ws.Range("N84").ClearContents()
$ws.Range("H113").Value = 762.0968
$ws.Range("I113").Value = 248.38889
$ws.Range("J113").Value = 1473.3846
$ws.Range("K113").Value = 745.1666700000001
$ws.Range("L113").Value = 4420.1538
$ws.Range("M113").Value = 1424.83333
$ws.Range("N113").Value = -8760.1538
